# Update paises.xlsx: reorder a handful of countries (which, due to the way
# the shared-string table is rebuilt, changes which row shows which country
# name) and refresh the case counts for several countries, plus the
# "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the "Datos actualizados" timestamp cell (row 1).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 10:59"

# ---------------------------------------------------------------------
# 2. Swap country-name pairs (rows keep their statistic columns tied to
#    the country whose numbers are listed there; where the underlying
#    data also changed, the new numbers are applied afterwards).
# ---------------------------------------------------------------------
$ws.Range("A32").Value  = "Indonesia"
$ws.Range("A33").Value  = "Emiratos Arabes Unidos"

$ws.Range("A165").Value = "Siria"
$ws.Range("A166").Value = "Islas Caimanes"

$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

# ---------------------------------------------------------------------
# 3. Refresh the statistic values (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#    rows whose numbers changed.
# ---------------------------------------------------------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2297360
$ws.Range("C4").Value = 170
$ws.Range("E4").Value = 1219876

# Rusia (row 6)
$ws.Range("B6").Value = 576952
$ws.Range("C6").Value = 7889
$ws.Range("D6").Value = 334592
$ws.Range("E6").Value = 234358
$ws.Range("G6").Value = 161
$ws.Range("H6").Value = 8002

# India (row 7)
$ws.Range("B7").Value = 396661
$ws.Range("C7").Value = 849
$ws.Range("D7").Value = 214346
$ws.Range("E7").Value = 169344
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 12971

# Banglades (row 20)
$ws.Range("B20").Value = 108775
$ws.Range("C20").Value = 3240
$ws.Range("D20").Value = 43993
$ws.Range("E20").Value = 63357
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = 1425

# Row 32 - now Indonesia (new figures)
$ws.Range("B32").Value = 45029
$ws.Range("C32").Value = 1226
$ws.Range("D32").Value = 17883
$ws.Range("E32").Value = 24717
$ws.Range("G32").Value = 56
$ws.Range("H32").Value = 2429

# Row 33 - now Emiratos Arabes Unidos (keeps its previous figures)
$ws.Range("B33").Value = 44145
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 30996
$ws.Range("E33").Value = 12849
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 300

# Polonia (row 39)
$ws.Range("B39").Value = 31620
$ws.Range("C39").Value = 304
$ws.Range("D39").Value = 16181
$ws.Range("E39").Value = 14093
$ws.Range("G39").Value = 12
$ws.Range("H39").Value = 1346

# Afganistan (row 42)
$ws.Range("B42").Value = 28424
$ws.Range("C42").Value = 546
$ws.Range("D42").Value = 8292
$ws.Range("E42").Value = 19563
$ws.Range("G42").Value = 21
$ws.Range("H42").Value = 569

# Moldavia (row 57)
$ws.Range("D57").Value = 7745
$ws.Range("E57").Value = 5353
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 458

# Lituania (row 112)
$ws.Range("B112").Value = 1795
$ws.Range("C112").Value = 3
$ws.Range("D112").Value = 1470
$ws.Range("E112").Value = 249

# Uganda (row 137)
$ws.Range("B137").Value = 763
$ws.Range("C137").Value = 8
$ws.Range("E137").Value = 271

# Row 165 - now Siria (new figures)
$ws.Range("B165").Value = 198
$ws.Range("C165").Value = 11
$ws.Range("D165").Value = 83
$ws.Range("E165").Value = 108
$ws.Range("H165").Value = 7

# Row 166 - now Islas Caimanes (keeps its previous figures)
$ws.Range("B166").Value = 195
$ws.Range("D166").Value = 143
$ws.Range("E166").Value = 51
$ws.Range("H166").Value = 1

# Row 208 - now Islas Turcas y Caicos (keeps its previous figures)
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209 - now Santa Sede (keeps its previous figures)
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
